$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (H)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 413
$wsOff.Range("C2").Value = 293
$wsOff.Range("D2").Value = 119
$wsOff.Range("E2").Value = 65

# DEF sheet - row 2 (H)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 426
$wsDef.Range("C2").Value = 301
$wsDef.Range("D2").Value = 114
$wsDef.Range("E2").Value = 57
$wsDef.Range("F2").Value = 10
$wsDef.Range("G2").Value = 5
